# PV-272 WIP changes:
#  - Rename the "Task Name" column header (C1) on the PV-Test-01 sheet to "Name"
#  - Switch the active/selected tab back to PV-Test-01 (it was left on "Dummy")
#    and restore the selection on that sheet to C2

$wb = $excel.ActiveWorkbook

$wsData  = $wb.Worksheets.Item("PV-Test-01")
$wsDummy = $wb.Worksheets.Item("Dummy")

# Header rename: "Task Name" -> "Name"
$wsData.Range("C1").Value = "Name"

# Make PV-Test-01 the active sheet/tab again, with C2 selected
$wsData.Activate()
$wsData.Range("C2").Select()
